$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.098888
$ws.Range("H2").Value = 6.296664
$ws.Range("I2").Value = 0.1082453658858517
$ws.Range("J2").Value = 0.1082453658858517
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.1030763333333333
$ws.Range("N2").Value = 0.309229
$ws.Range("O2").Value = 0.01126512502660735
$ws.Range("P2").Value = 0.01126512502660735
$ws.Range("Q2").Value = 0.2163456791173333
$ws.Range("R2").Value = 1.947111112056
$ws.Range("S2").Value = 0.001219397580254977
$ws.Range("T2").Value = 0.001219397580254977

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.098888
$ws.Range("H3").Value = 6.296664
$ws.Range("I3").Value = 0.1082453658858517
$ws.Range("J3").Value = 0.1082453658858517
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 9.046962666666667
$ws.Range("N3").Value = 27.140888
$ws.Range("O3").Value = 0.9887348749733926
$ws.Range("P3").Value = 0.9887348749733927
$ws.Range("Q3").Value = 18.98856137751467
$ws.Range("R3").Value = 170.897052397632
$ws.Range("S3").Value = 0.1070259683055967
$ws.Range("T3").Value = 0.1070259683055967

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 15.87514366666667
$ws.Range("H4").Value = 47.625431
$ws.Range("I4").Value = 0.8187243600843848
$ws.Range("J4").Value = 0.8187243600843847
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.1030763333333333
$ws.Range("N4").Value = 0.309229
$ws.Range("O4").Value = 0.01126512502660735
$ws.Range("P4").Value = 0.01126512502660735
$ws.Range("Q4").Value = 1.636351600299889
$ws.Range("R4").Value = 14.727164402699
$ws.Range("S4").Value = 0.009223032278679693
$ws.Range("T4").Value = 0.009223032278679693

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 15.87514366666667
$ws.Range("H5").Value = 47.625431
$ws.Range("I5").Value = 0.8187243600843848
$ws.Range("J5").Value = 0.8187243600843847
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 9.046962666666667
$ws.Range("N5").Value = 27.140888
$ws.Range("O5").Value = 0.9887348749733926
$ws.Range("P5").Value = 0.9887348749733927
$ws.Range("Q5").Value = 143.6218320803031
$ws.Range("R5").Value = 1292.596488722728
$ws.Range("S5").Value = 0.809501327805705
$ws.Range("T5").Value = 0.809501327805705

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.416064
$ws.Range("H6").Value = 4.248192
$ws.Range("I6").Value = 0.07303027402976368
$ws.Range("J6").Value = 0.07303027402976367
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.1030763333333333
$ws.Range("N6").Value = 0.309229
$ws.Range("O6").Value = 0.01126512502660735
$ws.Range("P6").Value = 0.01126512502660735
$ws.Range("Q6").Value = 0.1459626848853333
$ws.Range("R6").Value = 1.313664163968
$ws.Range("S6").Value = 0.0008226951676726839
$ws.Range("T6").Value = 0.0008226951676726838

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.416064
$ws.Range("H7").Value = 4.248192
$ws.Range("I7").Value = 0.07303027402976368
$ws.Range("J7").Value = 0.07303027402976367
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 9.046962666666667
$ws.Range("N7").Value = 27.140888
$ws.Range("O7").Value = 0.9887348749733926
$ws.Range("P7").Value = 0.9887348749733927
$ws.Range("Q7").Value = 12.81107814161066
$ws.Range("R7").Value = 115.299703274496
$ws.Range("S7").Value = 0.072207578862091
$ws.Range("T7").Value = 0.07220757886209099

